# correcting twitter files upload problem
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 already holds shared-string index 1 ("@wkupin" before the edit). The
# author reordered/retargeted the username list so that index now reads
# "@Ting2li" - update the text in place (style/format stays the General
# username style already applied to A2).
$ws.Range("A2").Value = "@Ting2li"

# A3 becomes a new username ("@MB_Leonard") and, like A4 below, switches
# to an explicit Text ("@") number format so the username string is never
# reinterpreted as a number/date by Excel.
$ws.Range("A3").Value = "@MB_Leonard"
$ws.Range("A3").NumberFormat = "@"

# A4 was empty before; it now holds the username that used to sit in A2
# ("@wkupin"), also with the Text number format applied.
$ws.Range("A4").Value = "@wkupin"
$ws.Range("A4").NumberFormat = "@"

# Move/save the active selection to D8, matching the author's last cursor
# position when the file was saved.
$ws.Range("D8").Select()
